# Apply updated cryptocurrency market data (price + 1h volume change)
# to the "cryptos" worksheet, matching the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "26.293.79"
$ws.Range("E2").Value = "  -0.03%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.690.55"
$ws.Range("E3").Value = "  +0.69%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.15%  "

# Row 5: BNB
$ws.Range("D5").Value = "'217.83"
$ws.Range("E5").Value = "  -0.03%  "

# Row 6: XRP
$ws.Range("D6").Value = "'0.5358"
$ws.Range("E6").Value = "  +1.87%  "

# Row 7: USDC
$ws.Range("E7").Value = "  -0.13%  "

# Row 8: Cardano
$ws.Range("E8").Value = "  +1.33%  "

# Row 9: Dogecoin
$ws.Range("E9").Value = "  -0.46%  "

# Row 10: Solana
$ws.Range("D10").Value = "'21.74"
$ws.Range("E10").Value = "  -0.85%  "

# Row 11: TRON
$ws.Range("D11").Value = "'0.07688"
$ws.Range("E11").Value = "  +2.48%  "

# Row 12: WrappedEther
$ws.Range("D12").Value = "1.712.05"
$ws.Range("E12").Value = "  +0.24%  "

# Row 13: Polkadot
$ws.Range("D13").Value = "'4.520"
$ws.Range("E13").Value = "  +0.18%  "

# Row 14: Polygon
$ws.Range("D14").Value = "'0.5799"
$ws.Range("E14").Value = "  +0.33%  "

# Row 15: ShibaInu
$ws.Range("D15").Value = "'0.000008372"
$ws.Range("E15").Value = "  -1.51%  "

# Row 16: Litecoin
$ws.Range("D16").Value = "'66.84"
$ws.Range("E16").Value = "  +3.19%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "26.347.09"
$ws.Range("E17").Value = "  +0.03%  "

# Row 18: Uniswap
$ws.Range("D18").Value = "'4.906"
$ws.Range("E18").Value = "  -0.28%  "

# Row 19: Dai
$ws.Range("E19").Value = "  -0.10%  "

# Row 20: Avalanche
$ws.Range("E20").Value = "  -0.06%  "

# Row 21: BitcoinCash
$ws.Range("D21").Value = "'193.88"
$ws.Range("E21").Value = "  +2.36%  "

# Row 22: Chainlink
$ws.Range("D22").Value = "'6.271"
$ws.Range("E22").Value = "  +1.23%  "

# Row 23: BinanceUSD
$ws.Range("E23").Value = "  -0.13%  "

# Row 24: Monero
$ws.Range("D24").Value = "'148.87"
$ws.Range("E24").Value = "  +2.78%  "

# Row 25: Stellar
$ws.Range("D25").Value = "'0.1286"
$ws.Range("E25").Value = "  +2.49%  "

# Row 26: Cosmos
$ws.Range("D26").Value = "'7.877"
$ws.Range("E26").Value = "  +1.43%  "

# Row 27: EthereumClassic
$ws.Range("E27").Value = "  +0.62%  "

# Row 28: Toncoin
$ws.Range("D28").Value = "'1.382"
$ws.Range("E28").Value = "  +1.40%  "

# Row 29: Hedera
$ws.Range("D29").Value = "'0.06120"
$ws.Range("E29").Value = "  -6.08%  "

# Row 30: PancakeSwap
$ws.Range("E30").Value = "  +0.17%  "

# Row 31: InternetComputer(DFINITY)
$ws.Range("D31").Value = "'3.602"
$ws.Range("E31").Value = "  +0.45%  "

# Row 32: Filecoin
$ws.Range("D32").Value = "'3.583"
$ws.Range("E32").Value = "  -0.23%  "

# Row 33: LidoDAOToken
$ws.Range("D33").Value = "'1.688"
$ws.Range("E33").Value = "  +1.76%  "

# Row 34: ARBITRUM
$ws.Range("E34").Value = "  +0.54%  "

# Row 35: ImmutableX
$ws.Range("D35").Value = "'0.6192"
$ws.Range("E35").Value = "  -0.18%  "

# Row 36: HuobiToken
$ws.Range("E36").Value = "  +0.85%  "

# Row 37: MXToken
$ws.Range("D37").Value = "'2.762"
$ws.Range("E37").Value = "  +0.75%  "

# Row 38: FraxShare
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "'6.211"
$ws.Range("E38").Value = "  -1.15%  "

# Row 39: VeChain
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01641"
$ws.Range("E39").Value = "  +1.25%  "

# Row 40: Maker
$ws.Range("D40").Value = "1.112.40"
$ws.Range("E40").Value = "  -0.42%  "

# Row 41: TrustWalletToken
$ws.Range("D41").Value = "'0.8779"
$ws.Range("E41").Value = "  +0.60%  "

# Row 42: PaxDollar
$ws.Range("E42").Value = "  -0.36%  "

# Row 43: Quant
$ws.Range("D43").Value = "'100.94"

# Row 44: RocketPoolETH
$ws.Range("D44").Value = "1.842.45"
$ws.Range("E44").Value = "  +0.70%  "

# Row 45: BabyDogeCoin
$ws.Range("D45").Value = "'0.00000000112"
$ws.Range("E45").Value = "  +4.27%  "

# Row 46: Aave
$ws.Range("D46").Value = "'57.80"
$ws.Range("E46").Value = "  +1.52%  "

# Row 47: Frax
$ws.Range("E47").Value = "  +0.37%  "

# Row 48: EnergySwap
$ws.Range("D48").Value = "'8.129"
$ws.Range("E48").Value = "  -0.26%  "

# Row 49: Cronos
$ws.Range("E49").Value = "  +0.39%  "

# Row 50: Mantle
$ws.Range("D50").Value = "'0.4292"
$ws.Range("E50").Value = "  -0.12%  "

# Row 51: Aptos
$ws.Range("E51").Value = "  -0.36%  "
